# "aggiunti dati misuraione parabole" - add the parabola measurement data
# to the "misure_meno_grossolane" sheet, rename the header in D1->A1's
# shared string from "m" to "d1(cm)", drop the old D1 header cell, and
# populate rows 2-19 with the new measurements (plus a handful of blank
# formatted rows below, 20-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("misure_meno_grossolane")

# --- header row -----------------------------------------------------
$ws.Cells.Item(1, 1).Value = "d1(cm)"
$ws.Cells.Item(1, 4).ClearContents()

# --- measurement data -------------------------------------------------
$data = @(
    @(5,  2.0211000000000001, 2.0550000000000002),
    @(10, 2.0024000000000002, 1.9923999999999999),
    @(15, 1.984,              1.9245000000000001),
    @(20, 1.9677,             1.8697999999999999),
    @(25, 1.9520999999999999, 1.8317000000000001),
    @(30, 1.9399,             1.8089),
    @(35, 1.9199999999999999, 1.7956000000000001),
    @(40, 1.9207000000000001, 1.7916000000000001),
    @(45, 1.9132,             1.7964),
    @(50, 1.9094,             1.8058000000000001),
    @(55, 1.9077,             1.8214999999999999),
    @(60, 1.9092,             1.8411),
    @(65, 1.9147000000000001, 1.8643000000000001),
    @(70, 1.9245000000000001, 1.8883000000000001),
    @(75, 1.9368000000000001, 1.9172),
    @(80, 1.9539,             1.9458),
    @(85, 1.9770000000000001, 1.9756),
    @(90, 2.0024000000000002, 2.0064000000000002)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 2).NumberFormat = "0.0000"
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 3).NumberFormat = "0.0000"
    $ws.Cells.Item($row, 4).NumberFormat = "0.0000"
    $row = $row + 1
}

# --- A14/A16/A18 pick up a distinct (but visually identical) style ----
foreach ($r in @(14, 16, 18)) {
    $ws.Cells.Item($r, 1).WrapText = $true
    $ws.Cells.Item($r, 1).WrapText = $false
}

# --- trailing blank but pre-formatted rows (20-23) --------------------
for ($r = 20; $r -le 23; $r++) {
    $ws.Cells.Item($r, 2).NumberFormat = "0.0000"
    $ws.Cells.Item($r, 3).NumberFormat = "0.0000"
    $ws.Cells.Item($r, 4).NumberFormat = "0.0000"
}

# --- row heights (default row height shrank to 14.25) -----------------
$ws.Rows.Item(1).RowHeight = 14.25
for ($r = 2; $r -le 23; $r++) {
    $ws.Rows.Item($r).RowHeight = 14.25
}
